$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.840.08'
$ws.Range('E2').Value = '  -0.38%  '
$ws.Range('D3').Value = '1.582.11'
$ws.Range('E3').Value = '  -2.40%  '
$ws.Range('E4').Value = '  -0.15%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '209.55'
$ws.Range('E5').Value = '  -1.40%  '
$ws.Range('E6').Value = '  -0.08%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.477'
$ws.Range('E7').Value = '  -1.96%  '
$ws.Range('E8').Value = '  -0.71%  '
$ws.Range('E9').Value = '  -1.09%  '
$ws.Range('E10').Value = '  -1.56%  '
$ws.Range('E11').Value = '  -0.29%  '
$ws.Range('D12').Value = '1.803.58'
$ws.Range('E12').Value = '  -2.31%  '
$ws.Range('D13').Value = '1.582.17'
$ws.Range('E13').Value = '  -2.29%  '
$ws.Range('E14').Value = '  -3.06%  '
$ws.Range('E15').Value = '  -2.55%  '
$ws.Range('D16').Value = '25.813.80'
$ws.Range('E16').Value = '  -0.54%  '
$ws.Range('D17').Value = '0.0₃0723'
$ws.Range('E17').Value = '  -2.11%  '
$ws.Range('E18').Value = '  -3.32%  '
$ws.Range('E19').Value = '  -0.09%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '191.26'
$ws.Range('E20').Value = '  -0.33%  '
$ws.Range('E21').Value = '  -2.05%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '9.31'
$ws.Range('E22').Value = '  -2.34%  '
$ws.Range('E23').Value = '  -1.29%  '
$ws.Range('E24').Value = '  -0.44%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '141.67'
$ws.Range('E25').Value = '  -1.85%  '
$ws.Range('E26').Value = '  -0.10%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '1.69'
$ws.Range('E27').Value = '  -0.68%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '15.05'
$ws.Range('E28').Value = '  -1.04%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '6.43'
$ws.Range('E29').Value = '  -3.18%  '
$ws.Range('E30').Value = '  -5.73%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '0.0470'
$ws.Range('E31').Value = '  -1.66%  '
$ws.Range('E32').Value = '  -0.25%  '
$ws.Range('E33').Value = '  -2.56%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '1.48'
$ws.Range('E34').Value = '  -0.51%  '
$ws.Range('D36').Value = '1.096.53'
$ws.Range('E36').Value = '  -2.82%  '
$ws.Range('E37').Value = '  -0.09%  '
$ws.Range('E38').Value = '  -2.24%  '
$ws.Range('E39').Value = '  -2.16%  '
$ws.Range('E40').Value = '  -3.41%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.774'
$ws.Range('E41').Value = '  -8.33%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.811'
$ws.Range('E42').Value = '  +7.48%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '93.71'
$ws.Range('E43').Value = '  -4.27%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '5.15'
$ws.Range('E44').Value = '  +0.15%  '
$ws.Range('D45').Value = '1.718.50'
$ws.Range('E45').Value = '  -2.22%  '
$ws.Range('D46').Value = '0.0₆0111'
$ws.Range('E46').Value = '  -1.43%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '1.49'
$ws.Range('E47').Value = '  -1.32%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '53.04'
$ws.Range('E49').Value = '  -1.60%  '
$ws.Range('E50').Value = '  -0.72%  '
$ws.Range('E51').Value = '  -0.11%  '
